# Update cryptocurrency price/volume data per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.132.43'
$ws.Range("D3").Value = '1.890.53'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''308.06'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").Value = '''0.9997'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '''0.5156'
$ws.Range("E7").Value = '  +2.41%  '
$ws.Range("D8").Value = '''0.3717'
$ws.Range("E8").Value = '  +1.84%  '
$ws.Range("D9").Value = '''0.07210'
$ws.Range("E9").Value = '  +0.67%  '
$ws.Range("D10").Value = '''0.9046'
$ws.Range("E10").Value = '  +1.50%  '
$ws.Range("D11").Value = '''21.02'
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("D12").Value = '''0.07624'
$ws.Range("E12").Value = '  +1.49%  '
$ws.Range("D13").Value = '1.893.02'
$ws.Range("E13").Value = '  +2.10%  '
$ws.Range("E14").Value = '  +3.14%  '
$ws.Range("D15").Value = '''5.277'
$ws.Range("E15").Value = '  +0.97%  '
$ws.Range("D16").Value = '''1.001'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '''0.000008511'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '''14.37'
$ws.Range("E18").Value = '  +2.20%  '
$ws.Range("D19").Value = '''0.9996'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = '27.162.24'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = '''5.057'
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("D22").Value = '2.133.36'
$ws.Range("E22").Value = '  +2.59%  '
$ws.Range("D23").Value = '''10.59'
$ws.Range("E23").Value = '  +2.59%  '
$ws.Range("D24").Value = '''6.432'
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").Value = '''145.30'
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("D27").Value = '''18.06'
$ws.Range("D28").Value = '''2.153'
$ws.Range("E28").Value = '  +4.47%  '
$ws.Range("D29").Value = '''114.64'
$ws.Range("E29").Value = '  +1.58%  '
$ws.Range("D30").Value = '''4.982'
$ws.Range("E30").Value = '  +6.95%  '
$ws.Range("D31").Value = '''4.808'
$ws.Range("E31").Value = '  +3.80%  '
$ws.Range("D32").Value = '''0.09211'
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("E33").Value = '  -0.58%  '
$ws.Range("D34").Value = '''1.199'
$ws.Range("E34").Value = '  +4.72%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''3.038'
$ws.Range("E35").Value = '  +2.28%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.7589'
$ws.Range("E36").Value = '  +2.81%  '
$ws.Range("D37").Value = '''3.273'
$ws.Range("E37").Value = '  +1.14%  '
$ws.Range("D38").Value = '''2.563'
$ws.Range("E38").Value = '  +1.93%  '
$ws.Range("D39").Value = '''0.5651'
$ws.Range("E39").Value = '  +6.28%  '
$ws.Range("D41").Value = '''1.077'
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '''8.965'
$ws.Range("E42").Value = '  +7.32%  '
$ws.Range("D43").Value = '''6.592'
$ws.Range("E43").Value = '  +1.81%  '
$ws.Range("D44").Value = '''118.26'
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("D45").Value = '''0.1508'
$ws.Range("E45").Value = '  +2.95%  '
$ws.Range("D46").Value = '''0.4819'
$ws.Range("E46").Value = '  +3.97%  '
$ws.Range("D47").Value = '''10.23'
$ws.Range("E47").Value = '  +3.02%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("E49").Value = '  +1.22%  '
$ws.Range("E50").Value = '  +0.76%  '
$ws.Range("D51").Value = '''63.57'
$ws.Range("E51").Value = '  +1.28%  '
